# actualizacion Vo.Bo. UPP 4T 2020 SIPOT
# Update the single-sheet "Reporte de Formatos" report: move the reporting
# period from 3T 2020 (Jul-Sep) to 4T 2020 (Oct-Dec), refresh the sign-off
# dates, reselect the header row, and let row heights re-fit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data row (row 8): period covered + review/sign-off dates ---------
# B8: Fecha de inicio del periodo que se informa  (2020-07-01 -> 2020-10-01)
$ws.Range("B8").Value = 44105
# C8: Fecha de término del periodo que se informa (2020-09-30 -> 2020-12-31)
$ws.Range("C8").Value = 44196
# K8 / L8: Fecha de validación / Fecha de actualización (2020-10-10 -> 2021-01-10)
$ws.Range("K8").Value = 44206
$ws.Range("L8").Value = 44206

# Row 8 no longer needs the manual 29.25pt height override - let Excel
# auto-fit it back to the default again.
$ws.Rows.Item(8).AutoFit()

# Row 3 (merged sub-header labels) now wraps onto two lines, so give it an
# explicit custom height.
$ws.Rows.Item(3).RowHeight = 27

# Move/restore the active selection to the top merged header cell (A2:C2)
# instead of leaving it parked on the data row.
$ws.Range("A2:C2").Select() | Out-Null
